$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by copying the existing "2022-Q1"
#    sheet (so it inherits the same header/number styles), placing it
#    right after "总计" and before "2022-Q1".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q1Sheet.Copy($null, $totalSheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# The copied sheet has only one data row (row 2); insert 4 more blank
# rows (rows 3-6) and stamp them with row 2's formatting so every row
# keeps the same styling as the template row.
$q4Sheet.Rows("3:6").Insert()
$q4Sheet.Rows("2").Copy()
$q4Sheet.Range("A3:H6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Populate the "2022-Q4" sheet with the new holdings data.
#    Numeric-looking values that must remain text (fund codes with
#    leading zeros, and decimal numbers stored as text like the other
#    quarter sheets) are written with a leading apostrophe so Excel
#    keeps them as text instead of silently converting to numbers.
# ---------------------------------------------------------------------
$q4Data = @(
    @("'004818", "国寿安保目标策略灵活配置混合A", "'2.76", "'59.92", "'2.87", "'0.0792", 7),
    @("'013067", "富安达中小盘六个月持有期混合",   "'2.09", "'78.87", "'2.91", "'0.0608", 10),
    @("'011629", "银河核心优势混合A",               "'2.29", "'69.64", "'2.55", "'0.0584", 10),
    @("'004819", "国寿安保目标策略灵活配置混合C", "'1.30", "'59.92", "'2.87", "'0.0373", 7),
    @("'016981", "银河核心优势混合C",               "'0.00", "'69.64", "'2.55", 0,        10)
)

for ($i = 0; $i -lt $q4Data.Count; $i++) {
    $r = $i + 2
    $row = $q4Data[$i]
    $q4Sheet.Cells.Item($r, 1).Value = $i
    $q4Sheet.Cells.Item($r, 2).Value = $row[0]
    $q4Sheet.Cells.Item($r, 3).Value = $row[1]
    $q4Sheet.Cells.Item($r, 4).Value = $row[2]
    $q4Sheet.Cells.Item($r, 5).Value = $row[3]
    $q4Sheet.Cells.Item($r, 6).Value = $row[4]
    $q4Sheet.Cells.Item($r, 7).Value = $row[5]
    $q4Sheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: insert a new row for "2022-Q4"
#    while keeping the existing "2022-Q1" / "2021-Q4" rows below it.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")

# Insert a blank row at row 3, pushing the existing "2021-Q4" row to row 4.
$ws1.Rows("3").Insert()

# Row 2 becomes the new "2022-Q4" entry.
$ws1.Cells.Item(2, 2).Value = "2022-Q4"
$ws1.Cells.Item(2, 3).Value = 5
$ws1.Cells.Item(2, 4).Value = 0.24

# Row 3 holds the data that used to be in row 2 ("2022-Q1").
$ws1.Cells.Item(3, 1).Value = 1
$ws1.Cells.Item(3, 2).Value = "2022-Q1"
$ws1.Cells.Item(3, 3).Value = 1
$ws1.Cells.Item(3, 4).Value = 0.35

# Match the formatting of the new A3 cell to its neighbours (A2/A4).
$ws1.Cells.Item(2, 1).Copy()
$ws1.Cells.Item(3, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 4 ("2021-Q4") keeps its original data; only its index changes.
$ws1.Cells.Item(4, 1).Value = 2

# Restore the original active sheet/selection (creating/copying sheets
# above shifts Excel's focus onto the last-touched sheet).
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
